# "Working grade lvl and floor restriction"
#
# 1) Sheet "7200": the locker 72027/1005641 pairing was removed (it was
#    being assigned by the wrong grade level / floor), so its row is
#    deleted and every row below it shifts up one. Dimension goes from
#    A1:B96 to A1:B95.
# 2) Sheet "Student Assignments": two new student/locker assignments are
#    appended (406527 -> 1005641, 412462 -> 1005641), extending the
#    dimension from A1:B6 to A1:B8.

$wb = $excel.ActiveWorkbook

# --- 1) Remove the stale locker row from the "7200" sheet ---
$ws7200 = $wb.Worksheets.Item("7200")
$ws7200.Rows(27).Delete()

# --- 2) Append the two new assignments on "Student Assignments" ---
# (stored as text, matching the existing ID#/Lockeruniq columns in this sheet)
$wsAssign = $wb.Worksheets.Item("Student Assignments")
$wsAssign.Range("A7:B8").NumberFormat = "@"
$wsAssign.Range("A7").Value = "406527"
$wsAssign.Range("B7").Value = "1005641"
$wsAssign.Range("A8").Value = "412462"
$wsAssign.Range("B8").Value = "1005641"
$wsAssign.Range("A7:B8").Style = "Normal"
